$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in A6 (continuing the sequence 1..5)
$ws.Range("A6").Value = 6

# Move the active selection to A7, as in the diff (selection activeCell="A7" sqref="A7")
$ws.Range("A7").Select()
